$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A146").Value = "Login with valid username and password"
$ws.Range("B146").Value = "PASSED"
$ws.Range("C146").Value = "edge"

$ws.Range("A147").Value = "Login with valid username and password"
$ws.Range("B147").Value = "PASSED"
$ws.Range("C147").Value = "edge"
